$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml) - column F = "想去人数" (interest count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 358
$ws1.Range("F4").Value = 1818
$ws1.Range("F9").Value = 156
$ws1.Range("F10").Value = 3588
$ws1.Range("F11").Value = 143
$ws1.Range("F12").Value = 95
$ws1.Range("F16").Value = 623
$ws1.Range("F17").Value = 117
$ws1.Range("F18").Value = 789
$ws1.Range("F19").Value = 11
$ws1.Range("F25").Value = 2834
$ws1.Range("F26").Value = 5293
$ws1.Range("F29").Value = 484
$ws1.Range("F30").Value = 3110
$ws1.Range("F31").Value = 302
$ws1.Range("F32").Value = 2297
$ws1.Range("F40").Value = 474
$ws1.Range("F41").Value = 820
$ws1.Range("F42").Value = 32
$ws1.Range("F44").Value = 456

# Sheet "演出" (sheet2.xml)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 76

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 358
$ws4.Range("F4").Value = 1818
$ws4.Range("F9").Value = 156
$ws4.Range("F10").Value = 3588
$ws4.Range("F11").Value = 143
$ws4.Range("F12").Value = 95
$ws4.Range("F14").Value = 76
$ws4.Range("F17").Value = 623
$ws4.Range("F18").Value = 117
$ws4.Range("F19").Value = 789
$ws4.Range("F20").Value = 11
$ws4.Range("F26").Value = 2834
$ws4.Range("F27").Value = 5293
$ws4.Range("F30").Value = 484
$ws4.Range("F31").Value = 0
$ws4.Range("F32").Value = 302
$ws4.Range("F33").Value = 2297
$ws4.Range("F41").Value = 474
$ws4.Range("F42").Value = 820
$ws4.Range("F43").Value = 32
